$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nlgn3"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2841056666666666
$ws.Range("H2").Value = 0.852317
$ws.Range("I2").Value = 0.1466007552634951
$ws.Range("J2").Value = 0.1466007552634951
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2191816666666667
$ws.Range("N2").Value = 0.657545
$ws.Range("O2").Value = 0.3167322406056181
$ws.Range("P2").Value = 0.3167322406056181
$ws.Range("Q2").Value = 0.06227075352944444
$ws.Range("R2").Value = 0.560436781765
$ws.Range("S2").Value = 0.04643318568908265
$ws.Range("T2").Value = 0.04643318568908266

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nlgn3"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2841056666666666
$ws.Range("H3").Value = 0.852317
$ws.Range("I3").Value = 0.1466007552634951
$ws.Range("J3").Value = 0.1466007552634951
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4728276666666666
$ws.Range("N3").Value = 1.418483
$ws.Range("O3").Value = 0.6832677593943819
$ws.Range("P3").Value = 0.6832677593943819
$ws.Range("Q3").Value = 0.1343330194567778
$ws.Range("R3").Value = 1.208997175111
$ws.Range("S3").Value = 0.1001675695744124
$ws.Range("T3").Value = 0.1001675695744124

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nlgn3"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8048609999999999
$ws.Range("H4").Value = 2.414583
$ws.Range("I4").Value = 0.4153145970881676
$ws.Range("J4").Value = 0.4153145970881677
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2191816666666667
$ws.Range("N4").Value = 0.657545
$ws.Range("O4").Value = 0.3167322406056181
$ws.Range("P4").Value = 0.3167322406056181
$ws.Range("Q4").Value = 0.176410775415
$ws.Range("R4").Value = 1.587696978735
$ws.Range("S4").Value = 0.1315435228919548
$ws.Range("T4").Value = 0.1315435228919548

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nlgn3"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8048609999999999
$ws.Range("H5").Value = 2.414583
$ws.Range("I5").Value = 0.4153145970881676
$ws.Range("J5").Value = 0.4153145970881677
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4728276666666666
$ws.Range("N5").Value = 1.418483
$ws.Range("O5").Value = 0.6832677593943819
$ws.Range("P5").Value = 0.6832677593943819
$ws.Range("Q5").Value = 0.380560548621
$ws.Range("R5").Value = 3.425044937589
$ws.Range("S5").Value = 0.2837710741962128
$ws.Range("T5").Value = 0.2837710741962128

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Nlgn3"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2762093333333333
$ws.Range("H6").Value = 0.8286279999999999
$ws.Range("I6").Value = 0.1425261852485395
$ws.Range("J6").Value = 0.1425261852485395
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2191816666666667
$ws.Range("N6").Value = 0.657545
$ws.Range("O6").Value = 0.3167322406056181
$ws.Range("P6").Value = 0.3167322406056181
$ws.Range("Q6").Value = 0.06054002202888889
$ws.Range("R6").Value = 0.54486019826
$ws.Range("S6").Value = 0.04514263799874129
$ws.Range("T6").Value = 0.04514263799874129

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Nlgn3"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2762093333333333
$ws.Range("H7").Value = 0.8286279999999999
$ws.Range("I7").Value = 0.1425261852485395
$ws.Range("J7").Value = 0.1425261852485395
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4728276666666666
$ws.Range("N7").Value = 1.418483
$ws.Range("O7").Value = 0.6832677593943819
$ws.Range("P7").Value = 0.6832677593943819
$ws.Range("Q7").Value = 0.1305994145915555
$ws.Range("R7").Value = 1.175394731324
$ws.Range("S7").Value = 0.09738354724979817
$ws.Range("T7").Value = 0.09738354724979817

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nlgn3"
$ws.Range("C8").Value = "Nrxn1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.572779
$ws.Range("H8").Value = 1.718337
$ws.Range("I8").Value = 0.2955584623997977
$ws.Range("J8").Value = 0.2955584623997977
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2191816666666667
$ws.Range("N8").Value = 0.657545
$ws.Range("O8").Value = 0.3167322406056181
$ws.Range("P8").Value = 0.3167322406056181
$ws.Range("Q8").Value = 0.1255426558516667
$ws.Range("R8").Value = 1.129883902665
$ws.Range("S8").Value = 0.09361289402583925
$ws.Range("T8").Value = 0.09361289402583925

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nlgn3"
$ws.Range("C9").Value = "Nrxn1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.572779
$ws.Range("H9").Value = 1.718337
$ws.Range("I9").Value = 0.2955584623997977
$ws.Range("J9").Value = 0.2955584623997977
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4728276666666666
$ws.Range("N9").Value = 1.418483
$ws.Range("O9").Value = 0.6832677593943819
$ws.Range("P9").Value = 0.6832677593943819
$ws.Range("Q9").Value = 0.2708257580856667
$ws.Range("R9").Value = 2.437431822771
$ws.Range("S9").Value = 0.2019455683739584
$ws.Range("T9").Value = 0.2019455683739584

